$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9978598947440904
$ws.Range("D2").Value = 331.1864081193348
$ws.Range("F2").Value = 1116.321384706528
$ws.Range("H2").Value = 4064.596306719501
$ws.Range("J2").Value = 11.05672974915624
$ws.Range("O2").Value = 0.4051833396466904
$ws.Range("R2").Value = 0.2778443611920325
$ws.Range("U2").Value = 133.8474128727878

$ws.Range("C3").Value = 0.9978264841463003
$ws.Range("D3").Value = 331.1974973578458
$ws.Range("F3").Value = 1116.305508020533
$ws.Range("H3").Value = 4249.843214534011
$ws.Range("J3").Value = 9.905387025253564
$ws.Range("O3").Value = 0.4092018282907502
$ws.Range("R3").Value = 0.2893276500200683
$ws.Range("U3").Value = 122.9840596055835

$ws.Range("C4").Value = 0.9978100511013438
$ws.Range("D4").Value = 331.202951886348
$ws.Range("F4").Value = 1116.319010303758
$ws.Range("H4").Value = 3842.899729375626
$ws.Range("J4").Value = 10.20424976095746
$ws.Range("O4").Value = 0.4060280184160411
$ws.Range("R4").Value = 0.2818596778159019
$ws.Range("U4").Value = 119.5369411197572

$ws.Range("C5").Value = 0.9977826199353224
$ws.Range("D5").Value = 331.2120573598033
$ws.Range("E5").Value = 1447.539256517647
$ws.Range("F5").Value = 1116.327199157843
$ws.Range("G5").Value = 44440.94208149001
$ws.Range("H5").Value = 3811.741752579173
$ws.Range("I5").Value = 106.9637844526894
$ws.Range("J5").Value = 10.36403416891608
$ws.Range("K5").Value = 0.3914122303326464
$ws.Range("N5").Value = 0.4092346670372784
$ws.Range("O5").Value = 0.4082685332737799
$ws.Range("R5").Value = 0.2859555945166232
$ws.Range("U5").Value = 117.3278186216055

$ws.Range("C6").Value = 0.9977781961815494
$ws.Range("D6").Value = 331.2135258230289
$ws.Range("E6").Value = 1447.541011877824
$ws.Range("F6").Value = 1116.327486054795
$ws.Range("G6").Value = 45288.6685051182
$ws.Range("H6").Value = 3947.474072645499
$ws.Range("I6").Value = 110.8975286659461
$ws.Range("J6").Value = 10.91844426340259
$ws.Range("K6").Value = 0.3933158787116008
$ws.Range("N6").Value = 0.4163762871582438
$ws.Range("O6").Value = 0.4098355853572959
$ws.Range("R6").Value = 0.2926683754256915
$ws.Range("U6").Value = 121.8159729293486

$ws.Range("C7").Value = 0.9977757378916616
$ws.Range("D7").Value = 331.2143418569635
$ws.Range("E7").Value = 1447.539632678328
$ws.Range("F7").Value = 1116.325290821365
$ws.Range("G7").Value = 46370.34152536967
$ws.Range("H7").Value = 4254.114625543581
$ws.Range("I7").Value = 116.1392125015711
$ws.Range("J7").Value = 11.14466441693333
$ws.Range("K7").Value = 0.3925913689007903
$ws.Range("N7").Value = 0.4170039008286233
$ws.Range("O7").Value = 0.4046797235506435
$ws.Range("R7").Value = 0.3126677993268265
$ws.Range("U7").Value = 127.2838769185045

$ws.Range("C8").Value = 0.9977600521657477
$ws.Range("D8").Value = 331.2195488577588
$ws.Range("E8").Value = 1447.534557917312
$ws.Range("F8").Value = 1116.315009059553
$ws.Range("G8").Value = 44105.4950117788
$ws.Range("H8").Value = 3915.578860588823
$ws.Range("I8").Value = 102.6329447904083
$ws.Range("J8").Value = 11.22393389599875
$ws.Range("K8").Value = 0.3949879232733975
$ws.Range("N8").Value = 0.405804338418588
$ws.Range("O8").Value = 0.4070157179008635
$ws.Range("R8").Value = 0.3022475171886146
$ws.Range("U8").Value = 113.856878686407

$ws.Range("C9").Value = 0.9977557134267026
$ws.Range("D9").Value = 331.2209891654111
$ws.Range("F9").Value = 1116.295905265526
$ws.Range("H9").Value = 3835.073084646454
$ws.Range("J9").Value = 11.09513910960653
$ws.Range("O9").Value = 0.4073026887334389
$ws.Range("R9").Value = 0.3142061820215845
$ws.Range("U9").Value = 119.1879125621929

$ws.Range("C10").Value = 0.9977568843190179
$ws.Range("D10").Value = 331.2206004694101
$ws.Range("E10").Value = 1447.52473694658
$ws.Range("F10").Value = 1116.30413647717
$ws.Range("G10").Value = 43156.99292466421
$ws.Range("H10").Value = 3790.648631944272
$ws.Range("I10").Value = 103.7940398613961
$ws.Range("J10").Value = 11.52920948130193
$ws.Range("K10").Value = 0.3925621554710648
$ws.Range("N10").Value = 0.4168855768793007
$ws.Range("O10").Value = 0.4084576999097317
$ws.Range("R10").Value = 0.3084445883022974
$ws.Range("U10").Value = 115.323249342698

